$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68, shifting existing rows 68:105 down to 69:106
$ws.Rows(68).Insert()

# Fill in the new row 68 with the new record (same fixed columns as the rest of
# the table, new variable data for D/I/J/K/L/M/P)
$ws.Cells.Item(68, 1).Value = 1
$ws.Cells.Item(68, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(68, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(68, 4).Value = 44830
$ws.Cells.Item(68, 5).Value = 15
$ws.Cells.Item(68, 6).Value = 100112038
$ws.Cells.Item(68, 7).Value = "Cebollín baby"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 300
$ws.Cells.Item(68, 11).Value = 1000
$ws.Cells.Item(68, 12).Value = 1500
$ws.Cells.Item(68, 13).Value = 1250
$ws.Cells.Item(68, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(68, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(68, 16).Value = 625
$ws.Cells.Item(68, 17).Value = 2
$ws.Cells.Item(68, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date number format as the rest of
# column D
$ws.Cells.Item(68, 4).NumberFormat = $ws.Cells.Item(69, 4).NumberFormat
